# Auto-generated edit script applying numeric corrections to the Asura_Profits
# sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1331.9032
$ws.Range("I15").Value = 1331.9032
$ws.Range("K15").Value = 3995.7096
$ws.Range("M15").Value = -3826.7096
# Row 33
$ws.Range("H33").Value = 195.53334
$ws.Range("I33").Value = 195.53334
$ws.Range("K33").Value = 195.53334
$ws.Range("M33").Value = 33.46665999999999
# Row 100
$ws.Range("H100").Value = 2379.5
$ws.Range("I100").Value = 1758.3334
$ws.Range("K100").Value = 1758.3334
$ws.Range("M100").Value = -1217.3334
# Row 101
$ws.Range("H101").Value = 755.8889
$ws.Range("I101").Value = 755.8889
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2267.6667
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -645.6667000000002
$ws.Range("N101").Value = ""
# Row 129
$ws.Range("H129").Value = 975.86206
$ws.Range("I129").Value = 475.75
$ws.Range("J129").Value = 1106.326
$ws.Range("K129").Value = 1427.25
$ws.Range("L129").Value = 3318.978
$ws.Range("M129").Value = 3572.75
$ws.Range("N129").Value = -13318.978
# Row 132
$ws.Range("H132").Value = 1517.7887
$ws.Range("I132").Value = 1408.409
$ws.Range("J132").Value = 2961.6
$ws.Range("K132").Value = 4225.227000000001
$ws.Range("L132").Value = 8884.799999999999
$ws.Range("M132").Value = -1695.227000000001
$ws.Range("N132").Value = -13944.8
# Row 137
$ws.Range("H137").Value = 1921.279
$ws.Range("I137").Value = 1761.1111
$ws.Range("J137").Value = 2191.5625
$ws.Range("K137").Value = 5283.3333
$ws.Range("L137").Value = 6574.6875
$ws.Range("M137").Value = -2733.3333
$ws.Range("N137").Value = -11674.6875
# Row 141
$ws.Range("H141").Value = 6181.6665
$ws.Range("I141").Value = 3684.6875
$ws.Range("J141").Value = 14172
$ws.Range("K141").Value = 11054.0625
$ws.Range("L141").Value = 42516
$ws.Range("M141").Value = -5874.0625
$ws.Range("N141").Value = -52876

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1803.8
$ws.Range("I61").Value = 1453.1538
$ws.Range("J61").Value = 4083
$ws.Range("K61").Value = 1453.1538
$ws.Range("L61").Value = 4083
$ws.Range("M61").Value = -1241.1538
$ws.Range("N61").Value = -4507
# Row 74
$ws.Range("H74").Value = 1269.1111
$ws.Range("I74").Value = 1107.9259
$ws.Range("J74").Value = 1752.6666
$ws.Range("K74").Value = 1107.9259
$ws.Range("L74").Value = 1752.6666
$ws.Range("M74").Value = -233.9259
$ws.Range("N74").Value = -3500.6666
# Row 77
$ws.Range("H77").Value = 1269.1111
$ws.Range("I77").Value = 1107.9259
$ws.Range("J77").Value = 1752.6666
$ws.Range("K77").Value = 5539.6295
$ws.Range("L77").Value = 8763.333000000001
$ws.Range("M77").Value = -1171.6295
$ws.Range("N77").Value = -17499.333
# Row 97
$ws.Range("H97").Value = 545.9091
$ws.Range("I97").Value = 570.5
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 570.5
$ws.Range("L97").Value = 300
$ws.Range("M97").Value = -74.5
$ws.Range("N97").Value = -1292
# Row 106
$ws.Range("H106").Value = 45000
$ws.Range("J106").Value = 45000
$ws.Range("L106").Value = 45000
$ws.Range("N106").Value = -47524
# Row 107
$ws.Range("H107").Value = 32076
$ws.Range("J107").Value = 32076
$ws.Range("L107").Value = 32076
$ws.Range("N107").Value = -39756
# Row 112
$ws.Range("H112").Value = 23126.215
$ws.Range("J112").Value = 23126.215
$ws.Range("L112").Value = 23126.215
$ws.Range("N112").Value = -26080.215
# Row 132
$ws.Range("H132").Value = 875454.9399999999
$ws.Range("I132").Value = 1113121.8
$ws.Range("J132").Value = 19854.4
$ws.Range("K132").Value = 3339365.4
$ws.Range("L132").Value = 59563.2
$ws.Range("M132").Value = -3336835.4
$ws.Range("N132").Value = -64623.2
# Row 136
$ws.Range("H136").Value = 1803.8
$ws.Range("I136").Value = 1453.1538
$ws.Range("J136").Value = 4083
$ws.Range("K136").Value = 4359.4614
$ws.Range("L136").Value = 12249
$ws.Range("M136").Value = -1809.4614
$ws.Range("N136").Value = -17349

$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 20933.572
$ws.Range("J88").Value = 20933.572
$ws.Range("L88").Value = 20933.572
$ws.Range("N88").Value = -21745.572
# Row 91
$ws.Range("H91").Value = 20933.572
$ws.Range("J91").Value = 20933.572
$ws.Range("L91").Value = 20933.572
$ws.Range("N91").Value = -23741.572
# Row 132
$ws.Range("H132").Value = 76099.336
$ws.Range("J132").Value = 76099.336
$ws.Range("L132").Value = 76099.336
$ws.Range("N132").Value = -86219.336
# Row 133
$ws.Range("H133").Value = 73350
$ws.Range("J133").Value = 73350
$ws.Range("L133").Value = 73350
$ws.Range("N133").Value = -83470
# Row 134
$ws.Range("H134").Value = 402662.4
$ws.Range("I134").Value = 608186.4399999999
$ws.Range("J134").Value = 3703.9412
$ws.Range("K134").Value = 1824559.32
$ws.Range("L134").Value = 11111.8236
$ws.Range("M134").Value = -1822024.32
$ws.Range("N134").Value = -16181.8236
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2907.3696
$ws.Range("I31").Value = 2047.9615
$ws.Range("J31").Value = 4024.6
$ws.Range("K31").Value = 2047.9615
$ws.Range("L31").Value = 4024.6
$ws.Range("M31").Value = -1752.9615
$ws.Range("N31").Value = -4614.6
# Row 34
$ws.Range("H34").Value = 2907.3696
$ws.Range("I34").Value = 2047.9615
$ws.Range("J34").Value = 4024.6
$ws.Range("K34").Value = 2047.9615
$ws.Range("L34").Value = 4024.6
$ws.Range("M34").Value = -1845.9615
$ws.Range("N34").Value = -4428.6
# Row 38
$ws.Range("H38").Value = 6019
$ws.Range("I38").Value = 38
$ws.Range("J38").Value = 12000
$ws.Range("K38").Value = 38
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 339
$ws.Range("N38").Value = -12754
# Row 46
$ws.Range("H46").Value = 6019
$ws.Range("I46").Value = 38
$ws.Range("J46").Value = 12000
$ws.Range("K46").Value = 38
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = 173
$ws.Range("N46").Value = -12422
# Row 58
$ws.Range("H58").Value = 1606.16
$ws.Range("I58").Value = 1472.7
$ws.Range("J58").Value = 2140
$ws.Range("K58").Value = 1472.7
$ws.Range("L58").Value = 2140
$ws.Range("M58").Value = -1269.7
$ws.Range("N58").Value = -2546
# Row 74
$ws.Range("H74").Value = 32754.666
$ws.Range("J74").Value = 32754.666
$ws.Range("L74").Value = 32754.666
$ws.Range("N74").Value = -34502.666
# Row 77
$ws.Range("H77").Value = 32754.666
$ws.Range("J77").Value = 32754.666
$ws.Range("L77").Value = 98263.99800000001
$ws.Range("N77").Value = -106999.998
# Row 94
$ws.Range("H94").Value = 1993.2142
$ws.Range("I94").Value = 1350
$ws.Range("K94").Value = 1350
$ws.Range("M94").Value = -899
# Row 132
$ws.Range("H132").Value = 2387.4285
$ws.Range("I132").Value = 1911.5294
$ws.Range("J132").Value = 3122.9092
$ws.Range("K132").Value = 5734.5882
$ws.Range("L132").Value = 9368.7276
$ws.Range("M132").Value = -3204.5882
$ws.Range("N132").Value = -14428.7276
# Row 134
$ws.Range("H134").Value = 1488.3055
$ws.Range("I134").Value = 1387.48
$ws.Range("J134").Value = 1717.4546
$ws.Range("K134").Value = 4162.440000000001
$ws.Range("L134").Value = 5152.3638
$ws.Range("M134").Value = -1627.440000000001
$ws.Range("N134").Value = -10222.3638
# Row 136
$ws.Range("H136").Value = 1606.16
$ws.Range("I136").Value = 1472.7
$ws.Range("J136").Value = 2140
$ws.Range("K136").Value = 4418.1
$ws.Range("L136").Value = 6420
$ws.Range("M136").Value = -1868.1
$ws.Range("N136").Value = -11520

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
# Row 132
$ws.Range("H132").Value = 1925.3334
$ws.Range("I132").Value = 1319.6471
$ws.Range("K132").Value = 3958.9413
$ws.Range("M132").Value = -1428.9413

$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 18000
$ws.Range("J64").Value = 18000
$ws.Range("L64").Value = 18000
$ws.Range("N64").Value = -18450
# Row 67
$ws.Range("H67").Value = 18000
$ws.Range("J67").Value = 18000
$ws.Range("L67").Value = 18000
$ws.Range("N67").Value = -19560
# Row 132
$ws.Range("H132").Value = 4187.604
$ws.Range("I132").Value = 4619.0293
$ws.Range("J132").Value = 3139.8572
$ws.Range("K132").Value = 13857.0879
$ws.Range("L132").Value = 9419.571599999999
$ws.Range("M132").Value = -11327.0879
$ws.Range("N132").Value = -14479.5716
# Row 136
$ws.Range("H136").Value = 2157.8914
$ws.Range("I136").Value = 2046.1945
$ws.Range("J136").Value = 2560
$ws.Range("K136").Value = 6138.583500000001
$ws.Range("L136").Value = 7680
$ws.Range("M136").Value = -3588.583500000001
$ws.Range("N136").Value = -12780

$ws = $wb.Worksheets.Item("WVR")
# Row 105
$ws.Range("H105").Value = 34500
$ws.Range("J105").Value = 34500
$ws.Range("L105").Value = 34500
$ws.Range("N105").Value = -41488
